# Applies the commit's language-neutralization + section-reorder edits.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. PROFESSIONAL SUMMARY paragraph: "all Black and Asian-American
#    voters" -> "50M voters"
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters, developed",
    2)

# ---------------------------------------------------------------------
# 2. "Partner - Siege Analytics" bullet: split the run so the new "50M"
#    token is bold + colored like the other emphasized numbers.
# ---------------------------------------------------------------------
$bullet = $d.Content
$null = $bullet.Find.Execute(
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$bulletStart = $bullet.Start
$bulletEnd = $bullet.End

# Replace the whole matched span with the plain-text target first so the
# surrounding runs keep their existing (unbolded) formatting...
$plain = $d.Range($bulletStart, $bulletEnd)
$plain.Text = "Discovered systematic race coding errors affecting 50M voters, developed"

# ...then re-find just the "50M" token within that span and bold/color it
# the same way the other highlighted metrics in this document are styled.
$numRange = $d.Content
$null = $numRange.Find.Execute("50M voters, developed geospatial machine learning", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tokenStart = $numRange.Start
$tokenRange = $d.Range($tokenStart, $tokenStart + 3)
$tokenRange.Font.Bold = 1
$tokenRange.Font.Color = 5258796

# ---------------------------------------------------------------------
# 3. Move the "Data Products Manager - Helm/Murmuration" block (5 paras)
#    from its old spot (after "Research Director - PCCC") to right after
#    the "Partner - Siege Analytics" block (i.e. before "Senior Analyst -
#    Myers Research").
# ---------------------------------------------------------------------
$r1 = $d.Content
$null = $r1.Find.Execute("Data Products Manager - Helm/Murmuration", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dpmStart = $r1.Paragraphs(1).Range.Start

$r2 = $d.Content
$null = $r2.Find.Execute("reducing processing time by", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dpmEnd = $r2.Paragraphs(1).Range.End

$dpmRange = $d.Range($dpmStart, $dpmEnd)
$dpmRange.Cut()

$r3 = $d.Content
$null = $r3.Find.Execute("Senior Analyst - Myers Research", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dpmTarget = $r3.Paragraphs(1).Range.Start
$dpmTargetRange = $d.Range($dpmTarget, $dpmTarget)
$dpmTargetRange.Paste()

# Pasting directly on a paragraph boundary makes the new first paragraph
# inherit the destination paragraph's style, so restore Heading 3 on it.
$dpmFixRange = $d.Range($dpmTarget, $dpmTarget)
$dpmFixPara = $dpmFixRange.Paragraphs(1)
$dpmFixPara.Style = "Heading 3"

# The Cut/Paste clipboard round-trip flattens multi-run paragraphs down to
# a single run, so the "57%" figure in the last bullet of this block loses
# its bold/colored emphasis. Re-apply it (only one "57%" in the document).
$pctRange = $d.Content
$null = $pctRange.Find.Execute("57%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pctStart = $pctRange.Start
$pctEnd = $pctRange.End
$pctFix = $d.Range($pctStart, $pctEnd)
$pctFix.Font.Bold = 1
$pctFix.Font.Color = 5258796

# ---------------------------------------------------------------------
# 4. Move the "Senior Analyst - Myers Research" block (5 paras) from its
#    old spot (right after "Partner - Siege Analytics") down to right
#    before "Research Director - PCCC" (i.e. after "Analytics
#    Supervisor - GSD&M"), which is where "Data Products Manager" used
#    to sit before step 3.
# ---------------------------------------------------------------------
$r4 = $d.Content
$null = $r4.Find.Execute("Senior Analyst - Myers Research", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$samStart = $r4.Paragraphs(1).Range.Start

$r5 = $d.Content
$null = $r5.Find.Execute("Co-developed a web application", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$samEnd = $r5.Paragraphs(1).Range.End

$samRange = $d.Range($samStart, $samEnd)
$samRange.Cut()

$r6 = $d.Content
$null = $r6.Find.Execute("Research Director - PCCC", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$samTarget = $r6.Paragraphs(1).Range.Start
$samTargetRange = $d.Range($samTarget, $samTarget)
$samTargetRange.Paste()

$samFixRange = $d.Range($samTarget, $samTarget)
$samFixPara = $samFixRange.Paragraphs(1)
$samFixPara.Style = "Heading 3"

# ---------------------------------------------------------------------
# 5. KEY PROJECTS / "Geospatial Demographic Classification System"
#    impact line: "all Black and Asian-American voters" -> "50M voters
#    nationwide"
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved",
    2)

Write-Output "All edits applied"
